# Localizes the remaining English boilerplate strings (and placeholder
# shape names) in the slide master, slide layouts, notes master and
# notes slide to Russian, matching the ru-RU "new.pptx" empty template.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Notes Master (ppt/notesMasters/notesMaster1.xml)
# ---------------------------------------------------------------------
$notesMaster = $p.NotesMaster

$notesMaster.Shapes.Item(1).Name = "Верхний колонтитул 1"   # id=2 Header Placeholder 1
$notesMaster.Shapes.Item(2).Name = "Дата 2"                  # id=3 Date Placeholder 2
$notesMaster.Shapes.Item(3).Name = "Рисунок 3"               # id=4 Slide Image Placeholder 3
$notesMaster.Shapes.Item(4).Name = "Заметка 4"               # id=5 Notes Placeholder 4
$notesMaster.Shapes.Item(5).Name = "Нижний колонтитул 5"     # id=6 Footer Placeholder 5
$notesMaster.Shapes.Item(6).Name = "Номер слайда 6"          # id=7 Slide Number Placeholder 6

$notesBody = $notesMaster.Shapes.Item(4).TextFrame.TextRange
$notesBody.Paragraphs(1, 1).Text = "Образец текста"
$notesBody.Paragraphs(2, 1).Text = "Второй уровень"
$notesBody.Paragraphs(3, 1).Text = "Третий уровень"
$notesBody.Paragraphs(4, 1).Text = "Четвертый уровень"
$notesBody.Paragraphs(5, 1).Text = "Пятый уровень"

# ---------------------------------------------------------------------
# 2. Notes Slide 1 (ppt/notesSlides/notesSlide1.xml)
# ---------------------------------------------------------------------
$notesSlide = $p.Slides.Item(1).NotesPage

$notesSlide.Shapes.Item(1).Name = "Рисунок 1"        # id=2 Slide Image Placeholder 1
$notesSlide.Shapes.Item(2).Name = "Текст 2"          # id=3 Notes Placeholder 2
$notesSlide.Shapes.Item(3).Name = "Номер слайда 3"   # id=4 Slide Number Placeholder 3

# ---------------------------------------------------------------------
# 3. Slide Master (ppt/slideMasters/slideMaster1.xml)
# ---------------------------------------------------------------------
$slideMaster = $p.SlideMaster

$slideMaster.Shapes.Item(1).Name = "Заголовок 1"                 # id=2 Title Placeholder 1
$slideMaster.Shapes.Item(2).Name = "Текст 2"                     # id=3 Text Placeholder 2
$slideMaster.Shapes.Item(3).Name = "Дата 3"                      # id=4 Date Placeholder 3
$slideMaster.Shapes.Item(4).Name = "Нижний колонтитул 4"         # id=5 Footer Placeholder 4
$slideMaster.Shapes.Item(5).Name = "Номер слайда 5"              # id=6 Slide Number Placeholder 5

$slideMaster.Shapes.Item(1).TextFrame.TextRange.Text = "Образец заголовка"

$masterBody = $slideMaster.Shapes.Item(2).TextFrame.TextRange
$masterBody.Paragraphs(1, 1).Text = "Образец текста"
$masterBody.Paragraphs(2, 1).Text = "Второй уровень"
$masterBody.Paragraphs(3, 1).Text = "Третий уровень"
$masterBody.Paragraphs(4, 1).Text = "Четвертый уровень"
$masterBody.Paragraphs(5, 1).Text = "Пятый уровень"

# ---------------------------------------------------------------------
# 4. Slide Layout "Заголовок и объект" (ppt/slideLayouts/slideLayout2.xml)
# ---------------------------------------------------------------------
$layoutObj = $slideMaster.CustomLayouts.Item(2)
$layoutObj.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5, 1).Text = "Пятый уровень"

# ---------------------------------------------------------------------
# 5. Slide Layout "Вертикальный заголовок и текст" (ppt/slideLayouts/slideLayout11.xml)
# ---------------------------------------------------------------------
$layoutVert = $slideMaster.CustomLayouts.Item(11)
$layoutVert.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5, 1).Text = "Пятый уровень"

# ---------------------------------------------------------------------
# 6. Slide Layout "Сравнение" (ppt/slideLayouts/slideLayout5.xml)
# ---------------------------------------------------------------------
$layoutCompare = $slideMaster.CustomLayouts.Item(5)
$layoutCompare.Shapes.Item(3).TextFrame.TextRange.Paragraphs(1, 1).Text = "Образец текста"

# ---------------------------------------------------------------------
# 7. Slide Layout "Рисунок с подписью" (ppt/slideLayouts/slideLayout9.xml)
# ---------------------------------------------------------------------
$layoutPic = $slideMaster.CustomLayouts.Item(9)
$layoutPic.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1, 1).Text = "Нажмите, чтобы добавить изображение"
